$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "blue - pink" color block rows (56 to 109 inclusive, 54 rows)
$ws.Range("A56:A109").EntireRow.Delete()
